$d = $word.ActiveDocument

# Package wrapper used for precise OOXML injection via Range.InsertXML.
# InsertXML *replaces* the exact contents of the range it is invoked on,
# so we always target either a freshly-created empty paragraph or the
# paragraph we want to rewrite in place.
function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:body>' + $bodyXml + '</w:body>' +
      '</w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'
}

function New-ParaXml([string]$text) {
    if ([string]::IsNullOrEmpty($text)) {
        return '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>'
    }
    return '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>' + $text + '</w:t></w:r></w:p>'
}

# Text for the new paragraphs, in document order. $null means "no run",
# i.e. an empty paragraph that only carries the eastAsia font hint.
$newParaTexts = @(
    $null,
    "优点：",
    "减少耦合：可以独立地开发、测试、优化、使用、理解和修改",
    "减轻维护的负担：可以更容易被程序员理解，并且在调试的时候可以不影响其他模块",
    "有效地调节性能：可以通过剖析确定哪些模块影响了系统的性能",
    "提高软件的可重用性",
    "降低了构建大型系统的风险：即使整个系统不可用，但是这些独立的模块却有可能是可用的"
)

# Locate the paragraph that holds the _GoBack bookmark: it is the empty
# paragraph immediately following the long "封装" description paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Contains("利用抽象数据类型")) {
        $targetIndex = $i + 1
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the anchor paragraph (the one right after the encapsulation description)."
}

$anchorPara = $d.Paragraphs.Item($targetIndex)
$anchorRange = $anchorPara.Range

foreach ($text in $newParaTexts) {
    [void]$anchorRange.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($targetIndex)
    [void]$newPara.Range.InsertXML((New-PkgXml (New-ParaXml $text)))
    $targetIndex = $targetIndex + 1
    $anchorPara = $d.Paragraphs.Item($targetIndex)
    $anchorRange = $anchorPara.Range
}

# Re-write the bookmark paragraph itself so it also carries the eastAsia
# paragraph-mark font hint, while preserving the _GoBack bookmark.
$bookmarkParaXml = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
[void]$anchorRange.InsertXML((New-PkgXml $bookmarkParaXml))
